# Adding load balanced Fargate service Terraform
#
# Reposition/resize a handful of shapes to make room for a second and
# third Availability Zone column, and relabel the AZ rectangles
# accordingly.
#
# Note: Shape.Left/Top/Width/Height are expressed in points (1 pt =
# 12700 EMU) but the point literals below are chosen so that, after the
# host's internal float handling, they land exactly on the target EMU
# value from the OOXML diff (rather than the "naive" emu/12700.0
# value, which can truncate one EMU short).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "Rectangle 138" - nudge down slightly (target y = 2606822 EMU)
$rect138 = $s.Shapes.Item(2)
$rect138.Top = 205.26158142089844

# "Graphic 140" - nudge down slightly, keeping pace with Rectangle 138
# (target y = 3412570 EMU)
$graphic140 = $s.Shapes.Item(3)
$graphic140.Top = 268.7063293457031

# "TextBox 139" ("Private subnet" label) - nudge down slightly
# (target y = 3467817 EMU)
$textbox139 = $s.Shapes.Item(10)
$textbox139.Top = 273.05645751953125

# "Rectangle 13" - the first Availability Zone box: grows taller/shifts
# up, and is relabeled "Availability Zone 2"
# (target x=4762347, y=307372, cx=2187552, cy=3639052 EMU)
$rect13 = $s.Shapes.Item(22)
$rect13.Left = 374.98797607421875
$rect13.Top = 24.2025203704834
$rect13.Width = 172.24819946289062
$rect13.Height = 286.53955078125
$rect13.TextFrame.TextRange.Text = "Availability Zone 2"

# "Rectangle 31" - the second Availability Zone box: grows taller/shifts
# up, and is relabeled "Availability Zone 3"
# (target x=7056413, y=307373, cx=2187552, cy=3655202 EMU)
$rect31 = $s.Shapes.Item(29)
$rect31.Left = 555.6231079101562
$rect31.Top = 24.202598571777344
$rect31.Width = 172.24819946289062
$rect31.Height = 287.8111877441406
$rect31.TextFrame.TextRange.Text = "Availability Zone 3"

# "Elbow Connector 77" - re-routed slightly to follow the resized shapes
# (target x=4444350, y=1194949, cx=529914, cy=2293832 EMU)
$elbow77 = $s.Shapes.Item(31)
$elbow77.Left = 349.9488220214844
$elbow77.Top = 94.0904769897461
$elbow77.Width = 41.72551345825195
$elbow77.Height = 180.61669921875

# "Elbow Connector 83" - adjust the bend point ratio (adj1 52195 -> 50478)
$elbow83 = $s.Shapes.Item(33)
$elbow83.Adjustments.Item(1) = 0.50478
